$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6583891
$ws.Range("I62").Value = 8932423
$ws.Range("K62").Value = 8932423
$ws.Range("M62").Value = -8931799
$ws.Range("H64").Value = 6201.633
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H65").Value = 6583891
$ws.Range("I65").Value = 8932423
$ws.Range("K65").Value = 44662115
$ws.Range("M65").Value = -44658995
$ws.Range("H67").Value = 6201.633
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H76").Value = 76929190
$ws.Range("I76").Value = 5497.6
$ws.Range("K76").Value = 5497.6
$ws.Range("M76").Value = -5182.6
$ws.Range("H79").Value = 76929190
$ws.Range("I79").Value = 5497.6
$ws.Range("K79").Value = 5497.6
$ws.Range("M79").Value = -4405.6
$ws.Range("H88").Value = 1320.5
$ws.Range("I88").Value = 653
$ws.Range("J88").Value = 1487.375
$ws.Range("K88").Value = 653
$ws.Range("L88").Value = 1487.375
$ws.Range("M88").Value = -247
$ws.Range("N88").Value = -2299.375
$ws.Range("H91").Value = 1320.5
$ws.Range("I91").Value = 653
$ws.Range("J91").Value = 1487.375
$ws.Range("K91").Value = 653
$ws.Range("L91").Value = 1487.375
$ws.Range("M91").Value = 751
$ws.Range("N91").Value = -4295.375
$ws.Range("H106").Value = 3424.05
$ws.Range("I106").Value = 3341.3684
$ws.Range("K106").Value = 3341.3684
$ws.Range("M106").Value = -2710.3684
$ws.Range("H125").Value = 7939906
$ws.Range("I125").Value = 4434
$ws.Range("J125").Value = 10104125
$ws.Range("K125").Value = 39906
$ws.Range("L125").Value = 90937125
$ws.Range("M125").Value = -37446
$ws.Range("N125").Value = -90942045
$ws.Range("H137").Value = 1737.9056
$ws.Range("I137").Value = 1591.0256
$ws.Range("J137").Value = 2147.0715
$ws.Range("K137").Value = 4773.0768
$ws.Range("L137").Value = 6441.2145
$ws.Range("M137").Value = -2223.0768
$ws.Range("N137").Value = -11541.2145
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5709.5347
$ws.Range("I32").Value = 5548.273
$ws.Range("K32").Value = 5548.273
$ws.Range("M32").Value = -5261.273
$ws.Range("H45").Value = 2978.2
$ws.Range("I45").Value = 2254.5715
$ws.Range("K45").Value = 2254.5715
$ws.Range("M45").Value = -1877.5715
$ws.Range("H97").Value = 1480.6875
$ws.Range("I97").Value = 1492.2858
$ws.Range("K97").Value = 1492.2858
$ws.Range("M97").Value = -996.2858000000001
$ws.Range("H102").Value = 1980.1333
$ws.Range("I102").Value = 1992.1818
$ws.Range("K102").Value = 1992.1818
$ws.Range("M102").Value = -370.1818000000001
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 149.27272
$ws.Range("I80").Value = 89.8
$ws.Range("K80").Value = 89.8
$ws.Range("M80").Value = 908.2
$ws.Range("H83").Value = 149.27272
$ws.Range("I83").Value = 89.8
$ws.Range("K83").Value = 449
$ws.Range("M83").Value = 4543
$ws.Range("H86").Value = 947833.75
$ws.Range("I86").Value = 1217557.4
$ws.Range("J86").Value = 3801
$ws.Range("K86").Value = 1217557.4
$ws.Range("L86").Value = 3801
$ws.Range("M86").Value = -1216434.4
$ws.Range("N86").Value = -6047
$ws.Range("H89").Value = 947833.75
$ws.Range("I89").Value = 1217557.4
$ws.Range("J89").Value = 3801
$ws.Range("K89").Value = 6087787
$ws.Range("L89").Value = 19005
$ws.Range("M89").Value = -6082171
$ws.Range("N89").Value = -30237
$ws.Range("H134").Value = 25404.066
$ws.Range("I134").Value = 3343.7778
$ws.Range("K134").Value = 10031.3334
$ws.Range("M134").Value = -7496.3334

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44344.668
$ws.Range("I31").Value = 1584.9
$ws.Range("J31").Value = 74887.36
$ws.Range("K31").Value = 1584.9
$ws.Range("L31").Value = 74887.36
$ws.Range("M31").Value = -1289.9
$ws.Range("N31").Value = -75477.36
$ws.Range("H34").Value = 44344.668
$ws.Range("I34").Value = 1584.9
$ws.Range("J34").Value = 74887.36
$ws.Range("K34").Value = 1584.9
$ws.Range("L34").Value = 74887.36
$ws.Range("M34").Value = -1382.9
$ws.Range("N34").Value = -75291.36
$ws.Range("H105").Value = 699.36365
$ws.Range("I105").Value = 623.6111
$ws.Range("K105").Value = 623.6111
$ws.Range("M105").Value = 1123.3889
$ws.Range("H134").Value = 252178.88
$ws.Range("I134").Value = 2199.2058
$ws.Range("K134").Value = 6597.617400000001
$ws.Range("M134").Value = -4062.617400000001
$ws.Range("H141").Value = 459487.22
$ws.Range("J141").Value = 565113.3
$ws.Range("L141").Value = 565113.3
$ws.Range("N141").Value = -575473.3

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 135.42105
$ws.Range("J2").Value = 166.61539
$ws.Range("L2").Value = 999.6923399999999
$ws.Range("N2").Value = -1225.69234
$ws.Range("H3").Value = 4333.3335
$ws.Range("J3").Value = 5000
$ws.Range("L3").Value = 15000
$ws.Range("N3").Value = -15224
$ws.Range("H44").Value = 349.66666
$ws.Range("I44").Value = 349.66666
$ws.Range("K44").Value = 1048.99998
$ws.Range("M44").Value = -650.9999800000001
$ws.Range("H56").Value = 6999.5
$ws.Range("I56").Value = 6999.5
$ws.Range("K56").Value = 6999.5
$ws.Range("M56").Value = -6469.5
$ws.Range("H92").Value = 1031.4375
$ws.Range("J92").Value = 1367.7142
$ws.Range("L92").Value = 4103.142599999999
$ws.Range("N92").Value = -6599.142599999999
$ws.Range("H97").Value = 462.85715
$ws.Range("I97").Value = 462.85715
$ws.Range("K97").Value = 1388.57145
$ws.Range("M97").Value = -892.5714499999999
$ws.Range("H109").Value = 49297.59
$ws.Range("I109").Value = 2108.1428
$ws.Range("J109").Value = 71319.336
$ws.Range("K109").Value = 6324.428400000001
$ws.Range("L109").Value = 213958.008
$ws.Range("M109").Value = -5284.428400000001
$ws.Range("N109").Value = -216038.008
$ws.Range("H121").Value = 835504.9399999999
$ws.Range("I121").Value = 2510
$ws.Range("J121").Value = 1113169.9
$ws.Range("K121").Value = 7530
$ws.Range("L121").Value = 3339509.7
$ws.Range("M121").Value = -6220
$ws.Range("N121").Value = -3342129.7
$ws.Range("H122").Value = 1981.826
$ws.Range("I122").Value = 1274.8182
$ws.Range("J122").Value = 2629.9167
$ws.Range("K122").Value = 11473.3638
$ws.Range("L122").Value = 23669.2503
$ws.Range("M122").Value = -9023.363799999999
$ws.Range("N122").Value = -28569.2503
$ws.Range("H133").Value = 20348.936
$ws.Range("J133").Value = 20807.334
$ws.Range("L133").Value = 62422.00199999999
$ws.Range("N133").Value = -72542.00199999999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2207.9111
$ws.Range("I102").Value = 1336.7812
$ws.Range("J102").Value = 4352.231
$ws.Range("K102").Value = 1336.7812
$ws.Range("L102").Value = 4352.231
$ws.Range("M102").Value = 285.2188000000001
$ws.Range("N102").Value = -7596.231
$ws.Range("H122").Value = 2454.8235
$ws.Range("I122").Value = 2264.3076
$ws.Range("K122").Value = 6792.9228
$ws.Range("M122").Value = -4342.9228
$ws.Range("H132").Value = 25302.09
$ws.Range("I132").Value = 3425.5676
$ws.Range("K132").Value = 10276.7028
$ws.Range("M132").Value = -7746.702799999999
$ws.Range("H136").Value = 29153.479
$ws.Range("J136").Value = 29153.479
$ws.Range("L136").Value = 87460.43700000001
$ws.Range("N136").Value = -92560.43700000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4895.7036
$ws.Range("I7").Value = 4526.722
$ws.Range("J7").Value = 5633.6665
$ws.Range("K7").Value = 4526.722
$ws.Range("L7").Value = 5633.6665
$ws.Range("M7").Value = -4414.722
$ws.Range("N7").Value = -5857.6665
$ws.Range("H126").Value = 4895.7036
$ws.Range("I126").Value = 4526.722
$ws.Range("J126").Value = 5633.6665
$ws.Range("K126").Value = 13580.166
$ws.Range("L126").Value = 16900.9995
$ws.Range("M126").Value = -11110.166
$ws.Range("N126").Value = -21840.9995
$ws.Range("H127").Value = 96564.336
$ws.Range("J127").Value = 96564.336
$ws.Range("L127").Value = 96564.336
$ws.Range("N127").Value = -106484.336

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6632.722
$ws.Range("I81").Value = 2324
$ws.Range("J81").Value = 15250.167
$ws.Range("K81").Value = 4648
$ws.Range("L81").Value = 30500.334
$ws.Range("M81").Value = -3587
$ws.Range("N81").Value = -32622.334
$ws.Range("H84").Value = 6632.722
$ws.Range("I84").Value = 2324
$ws.Range("J84").Value = 15250.167
$ws.Range("K84").Value = 23240
$ws.Range("L84").Value = 152501.67
$ws.Range("M84").Value = -17936
$ws.Range("N84").Value = -163109.67
$ws.Range("H122").Value = 41670504
$ws.Range("I122").Value = 62502932
$ws.Range("K122").Value = 187508796
$ws.Range("M122").Value = -187506346
$ws.Range("H132").Value = 15379.662
$ws.Range("I132").Value = 2262.6616
$ws.Range("J132").Value = 72220
$ws.Range("K132").Value = 6787.9848
$ws.Range("L132").Value = 216660
$ws.Range("M132").Value = -4257.9848
$ws.Range("N132").Value = -221720
$ws.Range("H136").Value = 58332.25
$ws.Range("I136").Value = 2769.0938
$ws.Range("K136").Value = 8307.2814
$ws.Range("M136").Value = -5757.2814
